$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=182; A='Medicine & Health'; B=6; C='Cardiology'; D='adversarial learning'; E='As a cardiologist, I want to utilize adversarial learning techniques to improve the robustness of my predictive models for heart disease risk assessment, ensuring that the models are resistant to adversarial attacks and provide reliable clinical insights.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=183; A='Medicine & Health'; B=6; C='Dermatology'; D='adversarial learning'; E='As a dermatologist, I want to leverage adversarial learning techniques to enhance the robustness of my skin condition classification models against potential adversarial attacks, ensuring reliable diagnoses and treatment recommendations.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=184; A='Medicine & Health'; B=6; C='Cardiology'; D='cnn'; E='As a researcher in cardiac imaging, I aim to develop CNN-based algorithms capable of detecting early signs of coronary artery disease from coronary CT angiography scans, enabling timely intervention and preventive care for patients at risk of heart attacks.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=185; A='Medicine & Health'; B=6; C='Dermatology'; D='cnn'; E='As a skincare researcher, I aim to develop CNN-based tools that can analyze facial photographs over time to track changes in skin conditions and assess the efficacy of skincare treatments, providing personalized recommendations for users.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=186; A='Medicine & Health'; B=6; C='Cardiology'; D='conversational agent'; E='As a patient, I want to interact with a conversational agent that uses machine learning to educate me about cardiovascular health, personalized risk factors, and lifestyle modifications, empowering me to make informed decisions for my well-being.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=187; A='Medicine & Health'; B=6; C='Dermatology'; D='conversational agent'; E='As a patient seeking dermatological advice, I want a conversational agent that uses machine learning to provide personalized skincare recommendations based on my skin type, concerns, and lifestyle habits, so that I can better manage and improve my skin health.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=188; A='Medicine & Health'; B=6; C='Cardiology'; D='decision tree'; E='As a healthcare provider, I want to use decision tree models to predict the likelihood of post-operative complications in cardiac surgery patients, based on pre-operative variables such as age, medical history, and surgical procedure type.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=189; A='Medicine & Health'; B=6; C='Dermatology'; D='decision tree'; E='As a dermatologist, I want to use decision tree algorithms to classify skin lesions based on visual characteristics such as size, color, and texture, so that I can make accurate diagnoses and recommend appropriate treatments.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=190; A='Medicine & Health'; B=6; C='Cardiology'; D='document classification'; E='As a cardiac rehabilitation specialist, I want a document classification tool that can analyze patient feedback forms and session reports to categorize patient progress and adherence to rehabilitation programs, helping to tailor personalized treatment plans for optimal recovery outcomes.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=191; A='Medicine & Health'; B=6; C='Dermatology'; D='document classification'; E='As a healthcare administrator, I need a document classification model specialized in dermatology to organize patient records, lab reports, and medical notes according to different skin conditions, ensuring efficient retrieval of information for clinical decision-making and patient management.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=192; A='Medicine & Health'; B=6; C='Cardiology'; D='entity extraction'; E='As a cardiologist, I need an entity extraction model in machine learning to automatically extract key clinical parameters such as blood pressure readings, cholesterol levels, and heart rate variations from patient medical records, facilitating comprehensive cardiovascular risk assessments.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=193; A='Medicine & Health'; B=6; C='Dermatology'; D='entity extraction'; E='As a skincare researcher, I need an entity extraction model capable of parsing clinical trial reports and extracting data on treatment efficacy metrics, adverse reactions, patient demographics, and treatment protocols for systematic analysis and comparison.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=194; A='Medicine & Health'; B=6; C='Cardiology'; D='feature selection'; E='As a cardiologist, I want to implement feature selection techniques in machine learning to identify the most significant biomarkers and clinical variables from cardiac imaging data, enhancing the accuracy of diagnostic models for detecting coronary artery disease.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=195; A='Medicine & Health'; B=6; C='Dermatology'; D='feature selection'; E='As a skincare researcher, I aim to use feature selection algorithms to analyze dermatological imaging data and identify key visual features (e.g., texture, color variation) associated with different types of skin lesions, aiding in automated diagnosis and treatment planning.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=196; A='Medicine & Health'; B=6; C='Cardiology'; D='imbalanced dataset'; E='As a pharmaceutical researcher, I want to develop algorithms using imbalanced dataset methodologies to evaluate the efficacy and safety of new cardiovascular drugs in clinical trials, ensuring robust analysis of rare adverse reactions and treatment outcomes.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=197; A='Medicine & Health'; B=6; C='Dermatology'; D='imbalanced dataset'; E='As a researcher studying dermatological adverse reactions, I aim to develop machine learning models that can effectively handle imbalanced datasets to predict the occurrence of rare but severe skin reactions to medications, facilitating early detection and prevention strategies.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=198; A='Medicine & Health'; B=6; C='Cardiology'; D='keyword extraction'; E='As a cardiologist, I want to implement keyword extraction algorithms in machine learning to automatically identify and extract key terms related to cardiac health from medical research papers, enabling efficient literature review and staying updated with the latest advancements in cardiology.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=199; A='Medicine & Health'; B=6; C='Dermatology'; D='keyword extraction'; E='As a dermatologist, I want to implement machine learning techniques for keyword extraction from patient medical histories and notes, to automatically identify and highlight key symptoms, treatments, and diagnostic indicators related to dermatological conditions, improving efficiency and accuracy in patient care.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=200; A='Medicine & Health'; B=6; C='Cardiology'; D='k-nearest neighbor'; E='As a cardiac rehabilitation specialist, I want to utilize k-Nearest Neighbor techniques to personalize exercise and rehabilitation plans for patients recovering from myocardial infarction, tailoring recommendations based on similarities to successful recovery cases.'; F='ReAdjusted_CoTPrompt'; EWrap=1 }
  @{ Row=201; A='Medicine & Health'; B=6; C='Dermatology'; D='k-nearest neighbor'; E='As a skincare researcher, I aim to apply k-Nearest Neighbor methods to analyze patient skincare routines and product usage patterns, identifying similarities among individuals with similar skin types and conditions to personalize skincare recommendations.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=202; A='Medicine & Health'; B=6; C='Cardiology'; D='multi-label classification'; E='As a cardiologist, I want to develop a multi-label classification model using machine learning to predict the presence of multiple cardiac conditions (e.g., hypertension, atrial fibrillation, coronary artery disease) from patient data, enabling comprehensive risk assessment and personalized treatment planning.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=203; A='Medicine & Health'; B=6; C='Dermatology'; D='multi-label classification'; E='As a dermatologist, I want to develop a multi-label classification model using machine learning to categorize skin conditions based on symptoms such as rash, itchiness, and discoloration, allowing for comprehensive diagnosis and treatment planning for patients with overlapping symptoms.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=204; A='Medicine & Health'; B=6; C='Cardiology'; D='neural network'; E='As a medical researcher, I aim to develop neural network models to predict patient outcomes following cardiac surgeries based on pre-operative risk factors, post-operative complications, and recovery progress, enhancing prognostic accuracy and patient care management.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=205; A='Medicine & Health'; B=6; C='Dermatology'; D='neural network'; E='As a skincare researcher, I aim to train a neural network for predictive modeling of skincare product effectiveness based on ingredients and user skin type data, allowing for personalized recommendations and improved consumer satisfaction.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=206; A='Medicine & Health'; B=6; C='Cardiology'; D='random forest'; E='As a cardiologist, I want to utilize a random forest algorithm to analyze patient data including age, lifestyle factors, and medical history to predict the likelihood of developing cardiovascular diseases such as coronary artery disease and heart failure, aiding in early intervention and preventive care.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=207; A='Medicine & Health'; B=6; C='Dermatology'; D='random forest'; E='As a dermatologist, I want to develop a random forest model to predict the likelihood of developing skin allergies based on patient demographics, environmental factors, and genetic predispositions, allowing for early preventive measures and personalized patient care.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=208; A='Medicine & Health'; B=6; C='Cardiology'; D='semantic similarity'; E='As a cardiologist, I want to develop a semantic similarity model in natural language processing to compare and categorize clinical notes and patient reports based on their relevance to specific cardiac conditions and symptoms, aiding in efficient information retrieval and decision-making.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=209; A='Medicine & Health'; B=6; C='Dermatology'; D='semantic similarity'; E='As a skincare product developer, I aim to use NLP techniques for semantic similarity to analyze customer reviews and feedback on skincare products, identifying common concerns and preferences among users to inform product improvement strategies.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=210; A='Medicine & Health'; B=6; C='Cardiology'; D='sentiment analysis'; E='As a healthcare provider, I want to apply sentiment analysis to patient feedback collected from cardiac rehabilitation programs, to understand patient satisfaction levels and identify areas for improvement in our services.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=211; A='Medicine & Health'; B=6; C='Dermatology'; D='sentiment analysis'; E='As a skincare product developer, I need to perform sentiment analysis on customer reviews and social media comments about our products to understand customer satisfaction levels and identify areas for product improvement in the skincare industry.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=212; A='Medicine & Health'; B=6; C='Cardiology'; D='speech to text'; E='As a cardiologist, I want to implement speech-to-text technology in clinical settings to transcribe patient interviews and discussions accurately, facilitating efficient documentation of symptoms, medical history, and treatment plans.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=213; A='Medicine & Health'; B=6; C='Dermatology'; D='speech to text'; E='As a dermatologist, I want to use speech to text technology during patient consultations to automatically transcribe discussions about symptoms, medical history, and treatment preferences, improving accuracy and efficiency in clinical documentation.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=214; A='Medicine & Health'; B=6; C='Cardiology'; D='text categorization'; E='As a healthcare provider, I need a text categorization system to classify medical research articles into categories such as coronary artery disease, heart failure, and arrhythmias, facilitating quick access to relevant literature for evidence-based practice in cardiology.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=215; A='Medicine & Health'; B=6; C='Dermatology'; D='text categorization'; E='As a skincare researcher, I want to implement text categorization algorithms to classify scientific articles and research papers in dermatology into topics such as skin cancer treatment, dermatological surgery techniques, and skincare product efficacy, facilitating literature review and research synthesis.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=216; A='Medicine & Health'; B=6; C='Cardiology'; D='unsupervised clustering'; E='As a cardiovascular researcher, I want to apply unsupervised clustering algorithms to cardiac imaging data (e.g., MRI, CT scans) to identify distinct patterns of heart morphology and function, aiding in the classification of structural heart diseases and anomalies.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=217; A='Medicine & Health'; B=6; C='Dermatology'; D='unsupervised clustering'; E='As a skincare product developer, I aim to use unsupervised clustering to analyze customer feedback and reviews on skincare products, identifying clusters of consumers with similar skin concerns and preferences to tailor product formulations and marketing strategies.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=218; A='Medicine & Health'; B=6; C='Cardiology'; D='voice recognition'; E='As a cardiac rehabilitation specialist, I aim to use voice recognition software for patients recovering from heart surgeries or cardiac events to record their daily progress and symptoms, facilitating remote monitoring and personalized care management.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=219; A='Medicine & Health'; B=6; C='Dermatology'; D='voice recognition'; E='As a dermatologist, I need a voice recognition system to transcribe patient consultations accurately, converting spoken descriptions of symptoms, medical history, and treatment preferences into text for efficient documentation and diagnosis.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=220; A='Medicine & Health'; B=6; C='Cardiology'; D='word embedding'; E='As a medical researcher in cardiology, I want to develop word embedding models to analyze and categorize medical literature and research articles on specific cardiovascular topics such as heart failure management and coronary artery disease prevention, enabling efficient literature review and knowledge synthesis.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
  @{ Row=221; A='Medicine & Health'; B=6; C='Dermatology'; D='word embedding'; E='As a dermatologist, I want to use word embedding techniques to analyze patient medical records and identify key dermatological terms and concepts related to symptoms, treatments, and disease progression, enhancing information retrieval and clinical decision-making.'; F='ReAdjusted_CoTPrompt'; EWrap=0 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $r.A
  $ws.Cells.Item($row, 2).Value = $r.B
  $ws.Cells.Item($row, 3).Value = $r.C
  $ws.Cells.Item($row, 4).Value = $r.D
  $ws.Cells.Item($row, 5).Value = $r.E
  $ws.Cells.Item($row, 6).Value = $r.F
  $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3)).Interior.Color = 15441517
  if ($r.EWrap -eq 1) {
    $ws.Cells.Item($row, 5).WrapText = $true
  }
}

$ws.Range("C196").Select()